$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "Karbonade m/ lauk"
$ws.Range("A3").Select()
